$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old data rows (content only, keep row/column formatting) so that the
# now-unused shared strings (FAPs, sCs, Fbn1, Itgb3, ECs) are dropped from the
# shared string table. They will be re-inserted below in the desired order.
$ws.Range("A2:T7").ClearContents()

# Sending / Target cluster labels, in the order they must first appear so the
# shared string table ends up ordered: ECs, FAPs, sCs, Fbn1, Itgb3
$clusters = @("ECs", "FAPs", "sCs")

$row = 2
$rowsInfo = @()
foreach ($sending in $clusters) {
    foreach ($target in $clusters) {
        $rowsInfo += , @($row, $sending, $target)
        $row++
    }
}

# Populate column A (Sending cluster) for every row first, then column D
# (Target cluster) for every row, so "ECs" ends up inserted into the shared
# string table before "FAPs" and "sCs".
foreach ($info in $rowsInfo) {
    $r = $info[0]
    $sending = $info[1]
    $ws.Cells.Item($r, 1).Value = $sending
}
foreach ($info in $rowsInfo) {
    $r = $info[0]
    $target = $info[2]
    $ws.Cells.Item($r, 4).Value = $target
}

# Now populate the constant Ligand / Receptor symbol columns (B, C).
foreach ($info in $rowsInfo) {
    $r = $info[0]
    $ws.Cells.Item($r, 2).Value = "Fbn1"
    $ws.Cells.Item($r, 3).Value = "Itgb3"
}

# Numeric statistics columns E..T for each row (values taken from the updated
# NATMI output).
$values = @{
    2  = @(2, 0.6666666666666666, 7.848425, 23.545275, 0.02436729568045431, 0.02436729568045431, 3, 1, 8.970048, 26.910144, 0.487108783009476, 0.4871087830094759, 70.4007489744, 633.6067407695999, 0.01186952374413816, 0.01186952374413816)
    3  = @(2, 0.6666666666666666, 7.848425, 23.545275, 0.02436729568045431, 0.02436729568045431, 3, 1, 9.012070666666666, 27.036212, 0.489390778604016, 0.489390778604016, 70.73056072203333, 636.5750464983, 0.01192512980553181, 0.01192512980553181)
    4  = @(2, 0.6666666666666666, 7.848425, 23.545275, 0.02436729568045431, 0.02436729568045431, 3, 1, 0.4327576666666667, 1.298273, 0.02350043838650813, 0.02350043838650813, 3.396466090008333, 30.568194810075, 0.0005726421307843421, 0.0005726421307843421)
    5  = @(3, 1, 291.329961, 873.989883, 0.9045029162236017, 0.9045029162236017, 3, 1, 8.970048, 26.910144, 0.487108783009476, 0.4871087830094759, 2613.243734008128, 23519.19360607315, 0.4405913147502006, 0.4405913147502006)
    6  = @(3, 1, 291.329961, 873.989883, 0.9045029162236017, 0.9045029162236017, 3, 1, 9.012070666666666, 27.036212, 0.489390778604016, 0.489390778604016, 2625.486195849244, 23629.37576264319, 0.4426553864202715, 0.4426553864202715)
    7  = @(3, 1, 291.329961, 873.989883, 0.9045029162236017, 0.9045029162236017, 3, 1, 0.4327576666666667, 1.298273, 0.02350043838650813, 0.02350043838650813, 126.075274152451, 1134.677467372059, 0.02125621505312967, 0.02125621505312967)
    8  = @(3, 1, 22.91008466666667, 68.730254, 0.07112978809594397, 0.07112978809594397, 3, 1, 8.970048, 26.910144, 0.487108783009476, 0.4871087830094759, 205.504559144064, 1849.541032296576, 0.03464794451513718, 0.03464794451513718)
    9  = @(3, 1, 22.91008466666667, 68.730254, 0.07112978809594397, 0.07112978809594397, 3, 1, 9.012070666666666, 27.036212, 0.489390778604016, 0.489390778604016, 206.4673019953164, 1858.205717957848, 0.03481026237821269, 0.03481026237821269)
    10 = @(3, 1, 22.91008466666667, 68.730254, 0.07112978809594397, 0.07112978809594397, 3, 1, 0.4327576666666667, 1.298273, 0.02350043838650813, 0.02350043838650813, 9.914514783482444, 89.230633051342, 0.00167158120259411, 0.00167158120259411)
}

foreach ($r in 2..10) {
    $rowValues = $values[$r]
    $col = 5
    foreach ($v in $rowValues) {
        $ws.Cells.Item($r, $col).Value = $v
        $col++
    }
}
